# parse treatments from tagging data, #809
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two stale header cell-comments (yellow-highlight / pp note) entirely.
$ws.Range("C1").Comment.Delete()
$ws.Range("N1").Comment.Delete()

# Rename a few headers (new shared strings, same cells) before the column
# layout shifts, while the old column letters still line up with the data.
$ws.Range("C1").Value = "PIT Tag #"
$ws.Range("O1").Value = "pp"
$ws.Range("X1").Value = "Amount"

# Drop the trailing, content-less header columns beyond "Universal Fish ID".
$ws.Range("AE1:AK1").EntireColumn.Delete()

# Drop the blank column that was merged into "double checking column ..."
$ws.Range("E1").EntireColumn.Delete()

# Re-point the conditional formatting rules that used to target the columns
# which have now shifted one to the left.
$condFmt1 = $ws.Range("J1").FormatConditions
$condFmt1.Item(1).ModifyAppliesToRange($ws.Range("I1"))
$condFmt1.Item(2).ModifyAppliesToRange($ws.Range("I1"))

$condFmt2 = $ws.Range("E1:E1048576").FormatConditions
$condFmt2.Item(1).ModifyAppliesToRange($ws.Range("D1:D1048576"))
$condFmt2.Item(2).ModifyAppliesToRange($ws.Range("D1:D1048576"))

# Rebuild the AutoFilter over the new, narrower header range.
$ws.AutoFilterMode = $false
$ws.Range("A1:AC1").AutoFilter()
